$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Bitcoin'
$ws.Range("B2").Value = 'btc'
$ws.Range("C2").Value = 98411
$ws.Range("D2").Value = 1950132842535
$ws.Range("E2").Value = 33488025597
$ws.Range("F2").Value = 1.24586

$ws.Range("A3").Value = 'Ethereum'
$ws.Range("B3").Value = 'eth'
$ws.Range("C3").Value = 2717.24
$ws.Range("D3").Value = 327567292210
$ws.Range("E3").Value = 17379584261
$ws.Range("F3").Value = 2.99203

$ws.Range("A4").Value = 'XRP'
$ws.Range("B4").Value = 'xrp'
$ws.Range("C4").Value = 2.52
$ws.Range("D4").Value = 145396047436
$ws.Range("E4").Value = 3895982591
$ws.Range("F4").Value = 4.81415

$ws.Range("A5").Value = 'Tether'
$ws.Range("B5").Value = 'usdt'
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 141898810750
$ws.Range("E5").Value = 64686637587
$ws.Range("F5").Value = 0.00393

$ws.Range("A6").Value = 'Solana'
$ws.Range("B6").Value = 'sol'
$ws.Range("C6").Value = 205.16
$ws.Range("D6").Value = 100083961971
$ws.Range("E6").Value = 4684942856
$ws.Range("F6").Value = 0.99521

$ws.Range("A7").Value = 'BNB'
$ws.Range("B7").Value = 'bnb'
$ws.Range("C7").Value = 637.47
$ws.Range("D7").Value = 92962391635
$ws.Range("E7").Value = 1194563777
$ws.Range("F7").Value = 4.67356

$ws.Range("A8").Value = 'USDC'
$ws.Range("B8").Value = 'usdc'
$ws.Range("C8").Value = 0.999996
$ws.Range("D8").Value = 56257179127
$ws.Range("E8").Value = 4785194982
$ws.Range("F8").Value = 0.00341

$ws.Range("A9").Value = 'Dogecoin'
$ws.Range("B9").Value = 'doge'
$ws.Range("C9").Value = 0.26772
$ws.Range("D9").Value = 39645674838
$ws.Range("E9").Value = 1351418423
$ws.Range("F9").Value = 6.91197

$ws.Range("A10").Value = 'Cardano'
$ws.Range("B10").Value = 'ada'
$ws.Range("C10").Value = 0.795423
$ws.Range("D10").Value = 28535539546
$ws.Range("E10").Value = 1157591917
$ws.Range("F10").Value = 14.58889

$ws.Range("A11").Value = 'Lido Staked Ether'
$ws.Range("B11").Value = 'steth'
$ws.Range("C11").Value = 2716.6
$ws.Range("D11").Value = 25553871100
$ws.Range("E11").Value = 62729785
$ws.Range("F11").Value = 3.07224

$ws.Range("A12").Value = 'TRON'
$ws.Range("B12").Value = 'trx'
$ws.Range("C12").Value = 0.246416
$ws.Range("D12").Value = 21201842858
$ws.Range("E12").Value = 829560532
$ws.Range("F12").Value = 5.871

$ws.Range("A13").Value = 'Wrapped Bitcoin'
$ws.Range("B13").Value = 'wbtc'
$ws.Range("C13").Value = 98186
$ws.Range("D13").Value = 12678178802
$ws.Range("E13").Value = 388622667
$ws.Range("F13").Value = 1.24356

$ws.Range("A14").Value = 'Chainlink'
$ws.Range("B14").Value = 'link'
$ws.Range("C14").Value = 19.69
$ws.Range("D14").Value = 12564995949
$ws.Range("E14").Value = 544695835
$ws.Range("F14").Value = 5.8177

$ws.Range("A15").Value = 'Wrapped stETH'
$ws.Range("B15").Value = 'wsteth'
$ws.Range("C15").Value = 3239.85
$ws.Range("D15").Value = 11043842730
$ws.Range("E15").Value = 58085274
$ws.Range("F15").Value = 3.09269

$ws.Range("A16").Value = 'Avalanche'
$ws.Range("B16").Value = 'avax'
$ws.Range("C16").Value = 26.69
$ws.Range("D16").Value = 10989021026
$ws.Range("E16").Value = 331081392
$ws.Range("F16").Value = 5.86355

$ws.Range("A17").Value = 'Sui'
$ws.Range("B17").Value = 'sui'
$ws.Range("C17").Value = 3.52
$ws.Range("D17").Value = 10870593863
$ws.Range("E17").Value = 1509548477
$ws.Range("F17").Value = 12.90979

$ws.Range("A18").Value = 'Stellar'
$ws.Range("B18").Value = 'xlm'
$ws.Range("C18").Value = 0.333881
$ws.Range("D18").Value = 10215378436
$ws.Range("E18").Value = 330524272
$ws.Range("F18").Value = 6.59628

$ws.Range("A19").Value = 'Litecoin'
$ws.Range("B19").Value = 'ltc'
$ws.Range("C19").Value = 129.46
$ws.Range("D19").Value = 9759833253
$ws.Range("E19").Value = 1918231053
$ws.Range("F19").Value = 14.06198

$ws.Range("A20").Value = 'Toncoin'
$ws.Range("B20").Value = 'ton'
$ws.Range("C20").Value = 3.87
$ws.Range("D20").Value = 9650081360
$ws.Range("E20").Value = 132078117
$ws.Range("F20").Value = 1.06107

$ws.Range("A21").Value = 'Shiba Inu'
$ws.Range("B21").Value = 'shib'
$ws.Range("C21").Value = 0.00001638
$ws.Range("D21").Value = 9641835463
$ws.Range("E21").Value = 211243492
$ws.Range("F21").Value = 3.16442

$ws.Range("A22").Value = 'Hedera'
$ws.Range("B22").Value = 'hbar'
$ws.Range("C22").Value = 0.246595
$ws.Range("D22").Value = 9434503192
$ws.Range("E22").Value = 291583721
$ws.Range("F22").Value = 3.76763

$ws.Range("A23").Value = 'LEO Token'
$ws.Range("B23").Value = 'leo'
$ws.Range("C23").Value = 9.9
$ws.Range("D23").Value = 9144770423
$ws.Range("E23").Value = 260436
$ws.Range("F23").Value = 0.03344

$ws.Range("A24").Value = 'USDS'
$ws.Range("B24").Value = 'usds'
$ws.Range("C24").Value = 1.001
$ws.Range("D24").Value = 8370764279
$ws.Range("E24").Value = 5957558
$ws.Range("F24").Value = 0.19937

$ws.Range("A25").Value = 'Hyperliquid'
$ws.Range("B25").Value = 'hype'
$ws.Range("C25").Value = 24.66
$ws.Range("D25").Value = 8243173619
$ws.Range("E25").Value = 115989134
$ws.Range("F25").Value = 4.31436

$ws.Range("A26").Value = 'WETH'
$ws.Range("B26").Value = 'weth'
$ws.Range("C26").Value = 2717.2
$ws.Range("D26").Value = 7937252033
$ws.Range("E26").Value = 885393927
$ws.Range("F26").Value = 3.00492

$ws.Range("A27").Value = 'Polkadot'
$ws.Range("B27").Value = 'dot'
$ws.Range("C27").Value = 5.16
$ws.Range("D27").Value = 7861278815
$ws.Range("E27").Value = 231265903
$ws.Range("F27").Value = 7.32178

$ws.Range("A28").Value = 'Bitget Token'
$ws.Range("B28").Value = 'bgb'
$ws.Range("C28").Value = 6.45
$ws.Range("D28").Value = 7743984697
$ws.Range("E28").Value = 279334002
$ws.Range("F28").Value = 0.80524

$ws.Range("A29").Value = 'Bitcoin Cash'
$ws.Range("B29").Value = 'bch'
$ws.Range("C29").Value = 341.87
$ws.Range("D29").Value = 6777263500
$ws.Range("E29").Value = 152775061
$ws.Range("F29").Value = 3.92785

$ws.Range("A30").Value = 'Ethena USDe'
$ws.Range("B30").Value = 'usde'
$ws.Range("C30").Value = 0.999782
$ws.Range("D30").Value = 6051454969
$ws.Range("E30").Value = 112005820
$ws.Range("F30").Value = -0.02404

$ws.Range("A31").Value = 'Uniswap'
$ws.Range("B31").Value = 'uni'
$ws.Range("C31").Value = 9.97
$ws.Range("D31").Value = 5994229870
$ws.Range("E31").Value = 299220810
$ws.Range("F31").Value = 9.1194

$ws.Range("A32").Value = 'MANTRA'
$ws.Range("B32").Value = 'om'
$ws.Range("C32").Value = 5.93
$ws.Range("D32").Value = 5766880575
$ws.Range("E32").Value = 222727811
$ws.Range("F32").Value = -3.07645

$ws.Range("A33").Value = 'Wrapped eETH'
$ws.Range("B33").Value = 'weeth'
$ws.Range("C33").Value = 2877.61
$ws.Range("D33").Value = 5471967311
$ws.Range("E33").Value = 15121965
$ws.Range("F33").Value = 3.04936

$ws.Range("A34").Value = 'Ondo'
$ws.Range("B34").Value = 'ondo'
$ws.Range("C34").Value = 1.4
$ws.Range("D34").Value = 4429222888
$ws.Range("E34").Value = 270640348
$ws.Range("F34").Value = 5.80717

$ws.Range("A35").Value = 'Pepe'
$ws.Range("B35").Value = 'pepe'
$ws.Range("C35").Value = 0.00001039
$ws.Range("D35").Value = 4365804086
$ws.Range("E35").Value = 852485531
$ws.Range("F35").Value = 8.92386

$ws.Range("A36").Value = 'Monero'
$ws.Range("B36").Value = 'xmr'
$ws.Range("C36").Value = 221.87
$ws.Range("D36").Value = 4086805671
$ws.Range("E36").Value = 64996379
$ws.Range("F36").Value = 0.24147

$ws.Range("A37").Value = 'WhiteBIT Coin'
$ws.Range("B37").Value = 'wbt'
$ws.Range("C37").Value = 27.73
$ws.Range("D37").Value = 3996396949
$ws.Range("E37").Value = 17794008
$ws.Range("F37").Value = 0.92908

$ws.Range("A38").Value = 'NEAR Protocol'
$ws.Range("B38").Value = 'near'
$ws.Range("C38").Value = 3.36
$ws.Range("D38").Value = 3957138556
$ws.Range("E38").Value = 227483967
$ws.Range("F38").Value = 4.92774

$ws.Range("A39").Value = 'Aave'
$ws.Range("B39").Value = 'aave'
$ws.Range("C39").Value = 260.45
$ws.Range("D39").Value = 3924874399
$ws.Range("E39").Value = 386194759
$ws.Range("F39").Value = 5.74941

$ws.Range("A40").Value = 'Aptos'
$ws.Range("B40").Value = 'apt'
$ws.Range("C40").Value = 6.31
$ws.Range("D40").Value = 3630205661
$ws.Range("E40").Value = 331319334
$ws.Range("F40").Value = 3.90311

$ws.Range("A41").Value = 'Mantle'
$ws.Range("B41").Value = 'mnt'
$ws.Range("C41").Value = 1.062
$ws.Range("D41").Value = 3569477708
$ws.Range("E41").Value = 106495921
$ws.Range("F41").Value = 3.20128

$ws.Range("A42").Value = 'Internet Computer'
$ws.Range("B42").Value = 'icp'
$ws.Range("C42").Value = 7.37
$ws.Range("D42").Value = 3538884457
$ws.Range("E42").Value = 90135187
$ws.Range("F42").Value = 4.82233

$ws.Range("A43").Value = 'Bittensor'
$ws.Range("B43").Value = 'tao'
$ws.Range("C43").Value = 431.42
$ws.Range("D43").Value = 3536199927
$ws.Range("E43").Value = 221947991
$ws.Range("F43").Value = 12.24481

$ws.Range("A44").Value = 'Dai'
$ws.Range("B44").Value = 'dai'
$ws.Range("C44").Value = 0.999786
$ws.Range("D44").Value = 3481684031
$ws.Range("E44").Value = 202816119
$ws.Range("F44").Value = -0.0143

$ws.Range("A45").Value = 'Official Trump'
$ws.Range("B45").Value = 'trump'
$ws.Range("C45").Value = 16.68
$ws.Range("D45").Value = 3334898510
$ws.Range("E45").Value = 927414180
$ws.Range("F45").Value = 4.89151

$ws.Range("A46").Value = 'sUSDS'
$ws.Range("B46").Value = 'susds'
$ws.Range("C46").Value = 1.038
$ws.Range("D46").Value = 3248904558
$ws.Range("E46").Value = 805276
$ws.Range("F46").Value = 0.43534

$ws.Range("A47").Value = 'Ethereum Classic'
$ws.Range("B47").Value = 'etc'
$ws.Range("C47").Value = 21.22
$ws.Range("D47").Value = 3198232671
$ws.Range("E47").Value = 98369155
$ws.Range("F47").Value = 4.18411

$ws.Range("A48").Value = 'OKB'
$ws.Range("B48").Value = 'okb'
$ws.Range("C48").Value = 50
$ws.Range("D48").Value = 2999139321
$ws.Range("E48").Value = 7582390
$ws.Range("F48").Value = 4.53236

$ws.Range("A49").Value = 'VeChain'
$ws.Range("B49").Value = 'vet'
$ws.Range("C49").Value = 0.03549073
$ws.Range("D49").Value = 2872965327
$ws.Range("E49").Value = 54254949
$ws.Range("F49").Value = 8.45239

$ws.Range("A50").Value = 'Gate'
$ws.Range("B50").Value = 'gt'
$ws.Range("C50").Value = 22.5
$ws.Range("D50").Value = 2837124409
$ws.Range("E50").Value = 15074510
$ws.Range("F50").Value = 4.45716

$ws.Range("A51").Value = 'POL (ex-MATIC)'
$ws.Range("B51").Value = 'pol'
$ws.Range("C51").Value = 0.327871
$ws.Range("D51").Value = 2802557458
$ws.Range("E51").Value = 77237287
$ws.Range("F51").Value = 6.57297
